# Update automatic: dades i banners [2026-02-07 09:19]
# Applies per-cell value updates to the active worksheet as described by the diff.
# Each cell is briefly switched to text format before assigning the new value so
# Excel does not auto-convert look-alike numbers/percentages/dates, then the
# number format is restored to General to match the original workbook styling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:42"
$r.NumberFormat = "General"

$r = $ws.Range("K2")
$r.NumberFormat = "@"
$r.Value = "0.3 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:44"
$r.NumberFormat = "General"

$r = $ws.Range("K3")
$r.NumberFormat = "@"
$r.Value = "0.4 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O3")
$r.NumberFormat = "@"
$r.Value = "-6.2 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:46"
$r.NumberFormat = "General"

$r = $ws.Range("J4")
$r.NumberFormat = "@"
$r.Value = "1002.0 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K4")
$r.NumberFormat = "@"
$r.Value = "0.7 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:49"
$r.NumberFormat = "General"

$r = $ws.Range("H5")
$r.NumberFormat = "@"
$r.Value = "76%"
$r.NumberFormat = "General"

$r = $ws.Range("J5")
$r.NumberFormat = "@"
$r.Value = "1002.0 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K5")
$r.NumberFormat = "@"
$r.Value = "0.5 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:51"
$r.NumberFormat = "General"

$r = $ws.Range("J6")
$r.NumberFormat = "@"
$r.Value = "1003.5 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K6")
$r.NumberFormat = "@"
$r.Value = "1.0 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O6")
$r.NumberFormat = "@"
$r.Value = "11.8 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:54"
$r.NumberFormat = "General"

$r = $ws.Range("H7")
$r.NumberFormat = "@"
$r.Value = "72%"
$r.NumberFormat = "General"

$r = $ws.Range("J7")
$r.NumberFormat = "@"
$r.Value = "1003.2 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K7")
$r.NumberFormat = "@"
$r.Value = "1.0 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:56"
$r.NumberFormat = "General"

$r = $ws.Range("H8")
$r.NumberFormat = "@"
$r.Value = "92%"
$r.NumberFormat = "General"

$r = $ws.Range("K8")
$r.NumberFormat = "@"
$r.Value = "1.1 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("L8")
$r.NumberFormat = "@"
$r.Value = "9.4 km/h - 190º 8:43 TU"
$r.NumberFormat = "General"

$r = $ws.Range("M8")
$r.NumberFormat = "@"
$r.Value = "11.5 °C 8:59 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O8")
$r.NumberFormat = "@"
$r.Value = "4.7 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:17:58"
$r.NumberFormat = "General"

$r = $ws.Range("O9")
$r.NumberFormat = "@"
$r.Value = "1.0 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:01"
$r.NumberFormat = "General"

$r = $ws.Range("H10")
$r.NumberFormat = "@"
$r.Value = "96%"
$r.NumberFormat = "General"

$r = $ws.Range("M10")
$r.NumberFormat = "@"
$r.Value = "12.5 °C 8:39 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O10")
$r.NumberFormat = "@"
$r.Value = "7.9 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:03"
$r.NumberFormat = "General"

$r = $ws.Range("J11")
$r.NumberFormat = "@"
$r.Value = "1006.0 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K11")
$r.NumberFormat = "@"
$r.Value = "0.3 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:05"
$r.NumberFormat = "General"

$r = $ws.Range("H12")
$r.NumberFormat = "@"
$r.Value = "66%"
$r.NumberFormat = "General"

$r = $ws.Range("K12")
$r.NumberFormat = "@"
$r.Value = "0.8 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("M12")
$r.NumberFormat = "@"
$r.Value = "13.5 °C 8:57 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O12")
$r.NumberFormat = "@"
$r.Value = "10.2 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:08"
$r.NumberFormat = "General"

$r = $ws.Range("H13")
$r.NumberFormat = "@"
$r.Value = "81%"
$r.NumberFormat = "General"

$r = $ws.Range("M13")
$r.NumberFormat = "@"
$r.Value = "13.8 °C 8:58 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O13")
$r.NumberFormat = "@"
$r.Value = "8.7 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:10"
$r.NumberFormat = "General"

$r = $ws.Range("K14")
$r.NumberFormat = "@"
$r.Value = "0.1 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:13"
$r.NumberFormat = "General"

$r = $ws.Range("J15")
$r.NumberFormat = "@"
$r.Value = "1002.3 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K15")
$r.NumberFormat = "@"
$r.Value = "0.6 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O15")
$r.NumberFormat = "@"
$r.Value = "5.7 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:15"
$r.NumberFormat = "General"

$r = $ws.Range("H16")
$r.NumberFormat = "@"
$r.Value = "93%"
$r.NumberFormat = "General"

$r = $ws.Range("K16")
$r.NumberFormat = "@"
$r.Value = "0.2 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:17"
$r.NumberFormat = "General"

$r = $ws.Range("J17")
$r.NumberFormat = "@"
$r.Value = "1005.5 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K17")
$r.NumberFormat = "@"
$r.Value = "0.3 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:20"
$r.NumberFormat = "General"

$r = $ws.Range("K18")
$r.NumberFormat = "@"
$r.Value = "0.5 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("M18")
$r.NumberFormat = "@"
$r.Value = "-5.8 °C 8:59 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O18")
$r.NumberFormat = "@"
$r.Value = "-7.7 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:22"
$r.NumberFormat = "General"

$r = $ws.Range("H19")
$r.NumberFormat = "@"
$r.Value = "99%"
$r.NumberFormat = "General"

$r = $ws.Range("J19")
$r.NumberFormat = "@"
$r.Value = "1006.9 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K19")
$r.NumberFormat = "@"
$r.Value = "0.9 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("L19")
$r.NumberFormat = "@"
$r.Value = "17.3 km/h - 298º 8:54 TU"
$r.NumberFormat = "General"

$r = $ws.Range("M19")
$r.NumberFormat = "@"
$r.Value = "6.3 °C 8:49 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O19")
$r.NumberFormat = "@"
$r.Value = "3.9 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:25"
$r.NumberFormat = "General"

$r = $ws.Range("K20")
$r.NumberFormat = "@"
$r.Value = "0.8 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("M20")
$r.NumberFormat = "@"
$r.Value = "-3.1 °C 8:38 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O20")
$r.NumberFormat = "@"
$r.Value = "-4.7 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:27"
$r.NumberFormat = "General"

$r = $ws.Range("J21")
$r.NumberFormat = "@"
$r.Value = "1002.7 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K21")
$r.NumberFormat = "@"
$r.Value = "0.5 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O21")
$r.NumberFormat = "@"
$r.Value = "5.6 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:30"
$r.NumberFormat = "General"

$r = $ws.Range("H22")
$r.NumberFormat = "@"
$r.Value = "86%"
$r.NumberFormat = "General"

$r = $ws.Range("K22")
$r.NumberFormat = "@"
$r.Value = "1.1 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("M22")
$r.NumberFormat = "@"
$r.Value = "12.9 °C 8:32 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O22")
$r.NumberFormat = "@"
$r.Value = "7.4 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:32"
$r.NumberFormat = "General"

$r = $ws.Range("H23")
$r.NumberFormat = "@"
$r.Value = "94%"
$r.NumberFormat = "General"

$r = $ws.Range("J23")
$r.NumberFormat = "@"
$r.Value = "1002.1 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K23")
$r.NumberFormat = "@"
$r.Value = "1.0 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("L23")
$r.NumberFormat = "@"
$r.Value = "25.2 km/h - 43º 8:59 TU"
$r.NumberFormat = "General"

$r = $ws.Range("M23")
$r.NumberFormat = "@"
$r.Value = "11.5 °C 8:59 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O23")
$r.NumberFormat = "@"
$r.Value = "7.8 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:35"
$r.NumberFormat = "General"

$r = $ws.Range("H24")
$r.NumberFormat = "@"
$r.Value = "78%"
$r.NumberFormat = "General"

$r = $ws.Range("J24")
$r.NumberFormat = "@"
$r.Value = "1001.3 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K24")
$r.NumberFormat = "@"
$r.Value = "1.0 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:37"
$r.NumberFormat = "General"

$r = $ws.Range("J25")
$r.NumberFormat = "@"
$r.Value = "1005.7 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:39"
$r.NumberFormat = "General"

$r = $ws.Range("H26")
$r.NumberFormat = "@"
$r.Value = "73%"
$r.NumberFormat = "General"

$r = $ws.Range("K26")
$r.NumberFormat = "@"
$r.Value = "0.2 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O26")
$r.NumberFormat = "@"
$r.Value = "-2.8 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:42"
$r.NumberFormat = "General"

$r = $ws.Range("H27")
$r.NumberFormat = "@"
$r.Value = "87%"
$r.NumberFormat = "General"

$r = $ws.Range("J27")
$r.NumberFormat = "@"
$r.Value = "1001.9 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K27")
$r.NumberFormat = "@"
$r.Value = "0.7 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O27")
$r.NumberFormat = "@"
$r.Value = "9.5 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:44"
$r.NumberFormat = "General"

$r = $ws.Range("H28")
$r.NumberFormat = "@"
$r.Value = "91%"
$r.NumberFormat = "General"

$r = $ws.Range("J28")
$r.NumberFormat = "@"
$r.Value = "1004.8 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:46"
$r.NumberFormat = "General"

$r = $ws.Range("K29")
$r.NumberFormat = "@"
$r.Value = "1.0 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("O29")
$r.NumberFormat = "@"
$r.Value = "10.6 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:49"
$r.NumberFormat = "General"

$r = $ws.Range("H30")
$r.NumberFormat = "@"
$r.Value = "76%"
$r.NumberFormat = "General"

$r = $ws.Range("K30")
$r.NumberFormat = "@"
$r.Value = "1.0 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("N30")
$r.NumberFormat = "@"
$r.Value = "-5.8 °C 8:50 TU"
$r.NumberFormat = "General"

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:51"
$r.NumberFormat = "General"

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:54"
$r.NumberFormat = "General"

$r = $ws.Range("H32")
$r.NumberFormat = "@"
$r.Value = "55%"
$r.NumberFormat = "General"

$r = $ws.Range("J32")
$r.NumberFormat = "@"
$r.Value = "1005.0 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K32")
$r.NumberFormat = "@"
$r.Value = "0.9 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("M32")
$r.NumberFormat = "@"
$r.Value = "12.8 °C 8:54 TU"
$r.NumberFormat = "General"

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:56"
$r.NumberFormat = "General"

$r = $ws.Range("M33")
$r.NumberFormat = "@"
$r.Value = "10.2 °C 8:48 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O33")
$r.NumberFormat = "@"
$r.Value = "7.1 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:18:58"
$r.NumberFormat = "General"

$r = $ws.Range("K34")
$r.NumberFormat = "@"
$r.Value = "0.3 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:19:01"
$r.NumberFormat = "General"

$r = $ws.Range("K35")
$r.NumberFormat = "@"
$r.Value = "0.2 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("N35")
$r.NumberFormat = "@"
$r.Value = "-8.9 °C 8:59 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O35")
$r.NumberFormat = "@"
$r.Value = "-6.4 °C"
$r.NumberFormat = "General"

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "2026-02-07 09:19:03"
$r.NumberFormat = "General"

$r = $ws.Range("J36")
$r.NumberFormat = "@"
$r.Value = "1007.3 hPa"
$r.NumberFormat = "General"

$r = $ws.Range("K36")
$r.NumberFormat = "@"
$r.Value = "0.9 MJ/m2"
$r.NumberFormat = "General"

$r = $ws.Range("M36")
$r.NumberFormat = "@"
$r.Value = "7.8 °C 8:59 TU"
$r.NumberFormat = "General"

$r = $ws.Range("O36")
$r.NumberFormat = "@"
$r.Value = "4.7 °C"
$r.NumberFormat = "General"

